$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 31 de Julio de 2020 a las 07:27"

# Row 16 - Pakistan
$ws.Range("B16").Value = 278305
$ws.Range("C16").Value = 903
$ws.Range("D16").Value = 247177
$ws.Range("E16").Value = 25177
$ws.Range("G16").Value = 27
$ws.Range("H16").Value = 5951

# Row 55 - Kirguistan
$ws.Range("B55").Value = 35619
$ws.Range("C55").Value = 396
$ws.Range("D55").Value = 24685
$ws.Range("E55").Value = 9562
$ws.Range("G55").Value = 8
$ws.Range("H55").Value = 1372

# Row 72 - Australia
$ws.Range("D72").Value = 9759
$ws.Range("E72").Value = 6948

# Row 110 - Tailandia
$ws.Range("B110").Value = 3310
$ws.Range("C110").Value = 6
$ws.Range("D110").Value = 3125
$ws.Range("E110").Value = 127

# Row 162 - Vietnam
$ws.Range("D162").Value = 373
$ws.Range("E162").Value = 136

# Row 192 - Papua Nueva Guinea
$ws.Range("B192").Value = 72
$ws.Range("C192").Value = 9
$ws.Range("E192").Value = 59

# Row 201 - Fiyi
$ws.Range("E201").Value = 8
$ws.Range("G201").Value = 1
$ws.Range("H201").Value = 1
